$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Sheet2"
$ws2.Range("B4").Value = "pepe"
$ws2.Range("B4").Select()

$ws1.Range("A1").Value = "lolo"
$ws1.Activate()
